# Invictus rigidité — "Cas de charge chassis left turn 2,2G"
#
# 1) Fill in the "Left Turn 2,2G" load-case results (D6:F15), which were
#    previously empty placeholders.
# 2) Swap the tab order of "Right Turn 2,2G" and "MAX SPEED" so that
#    "MAX SPEED" now sits right after "Left Turn 2,2G" and "Right Turn 2,2G"
#    moves to the very end (content is swapped between the two physical
#    sheets so each keeps its original tab position but trades data/name).
# 3) Make "Left Turn 2,2G" the active sheet/cell, since that is the case the
#    author just finished entering.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Left Turn 2,2G — fill the results table
# ---------------------------------------------------------------------
$leftTurn = $wb.Worksheets.Item("Left Turn 2,2G")

$leftTurnValues = @(
    @(0.261,  -2.39,  -2.9),
    @(-0.183, -2.39,  -2.15),
    @(0.258,  -1.91,  -2.4),
    @(-0.166, -1.92,  -1.78),
    @(0.243,  -1.51,  -1.89),
    @(-0.17,  -1.52,  -1.39),
    @(0.123,  -0.768, -0.662),
    @(-0.233, -0.77,  -0.686),
    @(0.123,  -0.494, -0.25),
    @(-0.138, -0.49,  -0.385)
)

for ($i = 0; $i -lt $leftTurnValues.Length; $i++) {
    $row = 6 + $i
    $vals = $leftTurnValues[$i]
    $leftTurn.Cells.Item($row, 4).Value = $vals[0]
    $leftTurn.Cells.Item($row, 5).Value = $vals[1]
    $leftTurn.Cells.Item($row, 6).Value = $vals[2]
}

# ---------------------------------------------------------------------
# 2) Swap "Right Turn 2,2G" <-> "MAX SPEED" (tab order)
# ---------------------------------------------------------------------
$rightTurn = $wb.Worksheets.Item("Right Turn 2,2G")
$maxSpeed  = $wb.Worksheets.Item("MAX SPEED")

# -- swap the C6:F15 data blocks --
for ($r = 6; $r -le 15; $r++) {
    for ($c = 4; $c -le 6; $c++) {
        $rtVal = $rightTurn.Cells.Item($r, $c).Value2
        $msVal = $maxSpeed.Cells.Item($r, $c).Value2
        $rightTurn.Cells.Item($r, $c).Value = $msVal
        $maxSpeed.Cells.Item($r, $c).Value  = $rtVal
    }
}

# -- swap the note cell: Right Turn keeps its note in I7, MAX SPEED's
#    note lives in J7, so trade both text and column --
$rightTurnNote = $rightTurn.Cells.Item(7, 9).Value2
$maxSpeedNote  = $maxSpeed.Cells.Item(7, 10).Value2

$rightTurn.Cells.Item(7, 9).ClearContents()
$maxSpeed.Cells.Item(7, 10).ClearContents()

$rightTurn.Cells.Item(7, 10).Value = $rightTurnNote
$maxSpeed.Cells.Item(7, 9).Value   = $maxSpeedNote

# -- swap the D2 title caption --
$rightTurn.Cells.Item(2, 4).Value = "MAX SPEED"
$maxSpeed.Cells.Item(2, 4).Value  = "Right Turn 2,2G"

# -- swap the remembered selection/scroll state of the two sheets BEFORE
#    renaming (so the $rightTurn / $maxSpeed handles below still line up
#    with "physically holds the Right Turn data" / "physically holds the
#    MAX SPEED data", regardless of what their .Name ends up being) --
# $rightTurn (index 6) takes over the data that used to live in MAX SPEED
# (index 7), so it should end up with MAX SPEED's old selection (E32), and
# vice-versa for $maxSpeed (index 7), which takes over Right Turn's old
# selection (H13).
$rightTurn.Activate()
$rightTurn.Range("E32").Select()

$maxSpeed.Activate()
$maxSpeed.Range("H13").Select()

# -- swap the sheet names themselves (via a temporary name, since both
#    names already exist in the workbook) --
$rightTurn.Name = "__TEMP_SWAP__"
$maxSpeed.Name  = "Right Turn 2,2G"
$rightTurn.Name = "MAX SPEED"

# ---------------------------------------------------------------------
# 3) Make "Left Turn 2,2G" the active sheet/cell
# ---------------------------------------------------------------------
$leftTurn.Activate()
$leftTurn.Range("H12").Select()
